$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column B (the "F/U" column), shifting column C ("NMYC") left into B
$ws.Range("B:B").Delete()
